$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 "Generator Data": add column C for "Generator 2" ---
$ws1.Range("C1").Value = "Generator 2"
$ws1.Range("B1").Copy()
$ws1.Range("C1").PasteSpecial(-4122)

$ws1.Range("A2").Value = "Nominal Capacity at upgrade 1"
$ws1.Range("B2").Value = 36091.293993
$ws1.Range("C2").Value = 1.65531644012

$ws1.Range("A3").Value = "Investment at upgrade 1"
$ws1.Range("B3").Value = 7218.2587986
$ws1.Range("C3").Value = 1.986379728144

$ws1.Range("A4").Value = "Yearly O&M Cost at upgrade 1"
$ws1.Range("B4").Value = 324.821645937
$ws1.Range("C4").Value = 0.08938708776647999

$ws1.Range("A5").Value = "Total actualized Fuel Cost"
$ws1.Range("B5").Value = 222674.729426
$ws1.Range("C5").Value = 5.17549102034

# --- Sheet2 "Yearly Fuel Costs": update existing values + extend to year 20 ---
$ws2.Range("B2").Value = 35575.04037179392
$ws2.Range("B3").Value = 35575.10697257982
$ws2.Range("B4").Value = 35575.17770987165
$ws2.Range("B5").Value = 35575.25858971571

$yearly = @(
    @(6, "Total Fuel Cost at y = 5", 35575.34992656275),
    @(7, "Total Fuel Cost at y = 6", 35575.45099812694),
    @(8, "Total Fuel Cost at y = 7", 35575.56455467194),
    @(9, "Total Fuel Cost at y = 8", 35575.69366701386),
    @(10, "Total Fuel Cost at y = 9", 35575.83646257433),
    @(11, "Total Fuel Cost at y = 10", 35575.99413476851),
    @(12, "Total Fuel Cost at y = 11", 35576.16809283406),
    @(13, "Total Fuel Cost at y = 12", 35576.36106154061),
    @(14, "Total Fuel Cost at y = 13", 35576.57596672527),
    @(15, "Total Fuel Cost at y = 14", 35576.81589628436),
    @(16, "Total Fuel Cost at y = 15", 35577.08357279628),
    @(17, "Total Fuel Cost at y = 16", 35577.38160892459),
    @(18, "Total Fuel Cost at y = 17", 35577.71326090474),
    @(19, "Total Fuel Cost at y = 18", 35578.08070180986),
    @(20, "Total Fuel Cost at y = 19", 35578.48787524875),
    @(21, "Total Fuel Cost at y = 20", 35578.93829007111)
)

foreach ($row in $yearly) {
    $r = $row[0]
    $label = $row[1]
    $value = $row[2]

    $ws2.Cells.Item($r, 1).Value = $label
    $ws2.Range("A2").Copy()
    $ws2.Cells.Item($r, 1).PasteSpecial(-4122)

    $ws2.Cells.Item($r, 2).Value = $value
}

Write-Output "done"
